$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows after the existing data (row 33), inheriting the
# formatting of the row directly above each new row (matches Excel's
# default "insert row" behaviour of carrying format down).
$ws.Rows.Item(34).Insert(-4121)
$ws.Rows.Item(35).Insert(-4121)
$ws.Rows.Item(36).Insert(-4121)

# Fill in the non-name/email columns first (these reuse existing shared
# strings so order doesn't matter for them).
$ws.Range("A34").Value = 110033
$ws.Range("B34").Value = 9317596771
$ws.Range("E34").Value = 818876434
$ws.Range("F34").Value = "ACT"
$ws.Range("G34").Value = "eng"
$ws.Range("H34").Value = "PWD"
$ws.Range("I34").Value = $true
$ws.Range("J34").Value = "superadmin"
$ws.Range("K34").Value = "now()"

$ws.Range("A35").Value = 110034
$ws.Range("B35").Value = 9317596772
$ws.Range("E35").Value = 818876435
$ws.Range("F35").Value = "ACT"
$ws.Range("G35").Value = "eng"
$ws.Range("H35").Value = "PWD"
$ws.Range("I35").Value = $true
$ws.Range("J35").Value = "superadmin"
$ws.Range("K35").Value = "now()"

$ws.Range("A36").Value = 110035
$ws.Range("B36").Value = 9317596773
$ws.Range("E36").Value = 818876436
$ws.Range("F36").Value = "ACT"
$ws.Range("G36").Value = "eng"
$ws.Range("H36").Value = "PWD"
$ws.Range("I36").Value = $true
$ws.Range("J36").Value = "superadmin"
$ws.Range("K36").Value = "now()"

# Names, then emails - entered column-by-column so new shared strings land
# in the same order as the source edit (all three names, then all three
# emails).
$ws.Range("C34").Value = "Nikola Tesla"
$ws.Range("C35").Value = "Graham Bell"
$ws.Range("C36").Value = "Albert Miles"

$ws.Range("D34").Value = "nikola.tesla@xyz.com"
$ws.Range("D35").Value = "graham.bell@xyz.com"
$ws.Range("D36").Value = "albert.miles@xyz.com"

# Select the row below the new data, matching the post-edit selection.
$ws.Range("A37:XFD1048576").Select()
